$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new user record as row 33, following the same column layout
# as the rest of the table:
#   A=id  B=uin  C=name  D=email  E=mobile  F=status_code  G=lang_code
#   H=last_login_method  I=is_active  J=cr_by  K=cr_dtimes
$row = 33
$ws.Cells.Item($row, 1).Value  = 110032
$ws.Cells.Item($row, 2).Value  = 9317596770
$ws.Cells.Item($row, 3).Value  = "Ewan Marsh"
$ws.Cells.Item($row, 4).Value  = "ewan.marsh@xyz.com"
$ws.Cells.Item($row, 5).Value  = 818876433
$ws.Cells.Item($row, 6).Value  = "ACT"
$ws.Cells.Item($row, 7).Value  = "eng"
$ws.Cells.Item($row, 8).Value  = "PWD"
$ws.Cells.Item($row, 9).Value  = $true
$ws.Cells.Item($row, 10).Value = "superadmin"
$ws.Cells.Item($row, 11).Value = "now()"

# Match the visible formatting used by the rest of the table: the
# "is_active" column is left-aligned.
$ws.Cells.Item($row, 9).HorizontalAlignment = -4131

# Reset the view: select the columns to the right of the data (mirrors
# a "select remaining columns" action) and scroll back to the top of
# the sheet instead of staying scrolled down near the previous last row.
$ws.Columns("L:XFD").Select()
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
